$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "RECOVERY" trial row (row 10) was removed from the Query1 table.
# Delete the entire worksheet row so everything below shifts up, the
# table/dimension ranges shrink by one row, and the now-unused
# "RECOVERY" shared string drops out on save.
$ws.Rows.Item(10).Delete()

# The hidden ExternalData_1 defined name tracks the query's cached
# extent and needs to be shrunk to match the smaller result set.
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$B`$15"
